$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-75 down to 12-76.
$ws.Rows.Item(11).Insert(-4121)

# Populate the newly inserted row 11 with the latest weekly price report.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(11, 3).Value = 'La Araucanía'
$ws.Cells.Item(11, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(11, 4).Value = 44613
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 100112030
$ws.Cells.Item(11, 7).Value = 'Poroto granado'
$ws.Cells.Item(11, 8).Value = 'Sin especificar'
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 70
$ws.Cells.Item(11, 11).Value = 25000
$ws.Cells.Item(11, 12).Value = 28000
$ws.Cells.Item(11, 13).Value = 26286
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(11, 16).Value = 1051
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = 'Hortaliza'
